$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws2 = $wb.Worksheets.Item("DATA")

# --- Sheet1 (RUNMANAGER) edits ---

# Row2: execute value changes from "yes" to "no"
$ws1.Range("C2").Value = "no"

# New row 4
$ws1.Range("A4").Value = "amazonTest"
$ws1.Range("B4").Value = "check"
$ws1.Range("C4").Value = "Yes"
$ws1.Range("D4").Value = "'1"
$ws1.Range("E4").Value = "'1"

# --- Sheet2 (DATA) edits ---

# New column G header + existing rows get quote-prefixed empty text in col G
$ws2.Range("G1").Value = "menutext"
$ws2.Range("G2").Value = "'"
$ws2.Range("G3").Value = "'"
$ws2.Range("G4").Value = "'"
$ws2.Range("G5").Value = "'"

# New row 6
$ws2.Range("A6").Value = "amazonTest"
$ws2.Range("B6").Value = "yes"
$ws2.Range("C6").Value = "chrome"
$ws2.Range("D6").Value = "'"
$ws2.Range("E6").Value = "'"
$ws2.Range("F6").Value = "'"
$ws2.Range("G6").Value = "Laptops"

# --- Selection / active sheet state ---
# RUNMANAGER selection moves to A4
$ws1.Range("A4").Select() | Out-Null
# DATA is the active/tabSelected sheet, selection moves to H7; select last so DATA stays active
$ws2.Range("H7").Select() | Out-Null
